# Update "想去人数" (want-to-go count) values on the 展览 (Exhibition) and
# 全部类型 (All Types) worksheets, per the upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 35
$wsExhibition.Range("F4").Value = 16144
$wsExhibition.Range("F8").Value = 15548
$wsExhibition.Range("F10").Value = 9178
$wsExhibition.Range("F28").Value = 513
$wsExhibition.Range("F32").Value = 74
$wsExhibition.Range("F39").Value = 5646

# --- Sheet "全部类型" (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 35
$wsAll.Range("F4").Value = 16144
$wsAll.Range("F8").Value = 15548
$wsAll.Range("F10").Value = 9178
$wsAll.Range("F28").Value = 513
$wsAll.Range("F34").Value = 74
$wsAll.Range("F41").Value = 5646
